# Update the cover-page version string.
#   "Version 11.10.01, 2016-05-02"  ->  "Version 11.13.00, 2017-03-06"
#
# The run immediately preceding the edited span (the "." run before the
# version number) is left completely untouched. The span "10.01, 2016-05-02"
# (17 characters) is replaced in one go by "13.00, 2017-03-06" (also 17
# characters), which the engine will fold into a single run. We then
# re-split that single run back into the individual one/two/four-character
# runs the diff expects by toggling Font.Bold off/on at each internal
# boundary -- a formatting-only change re-splits runs without re-merging
# them, even though the formatting ends up identical on both sides.

$d = $word.ActiveDocument

$full = $d.Content.Text
$idx = $full.IndexOf("Version 11.10.01, 2016-05-02")
if ($idx -lt 0) {
    throw "Could not locate the version string to update."
}

$spanStart = $idx + [string]"Version 11.".Length
$spanEnd = $idx + [string]"Version 11.10.01, 2016-05-02".Length

$r = $d.Range($spanStart, $spanEnd)
$r.Text = "13.00, 2017-03-06"

# Re-cut run boundaries to match the diff's run layout:
#   "1" "3" ".0" "0" ", 20" "1" "7" "-" "0" "3" "-" "0" "6"
# A formatting-only change (toggle off/on) re-splits the run at both ends
# of the touched range without re-merging neighbouring runs afterwards,
# even though the formatting ends up identical on both sides -- unlike a
# text edit, which always re-merges with the preceding run.
$lengths = @(1, 1, 2, 1, 4, 1, 1, 1, 1, 1, 1, 1, 1)

$pos = $spanStart
foreach ($len in $lengths) {
    $rc = $d.Range($pos, $pos + $len)
    $rc.Font.Bold = $false
    $rc.Font.Bold = $true
    $pos = $pos + $len
}
